$d = $word.ActiveDocument

# Grab a range collapsed to the very start of the document
$r = $d.Range(0, 0)

# Insert the two new heading paragraphs before the existing content
$r.InsertBefore("The chivalry of War`rPart 1`r")
